$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the old "_old"/"_new" suffixed header labels to the new
# "_FV2210"/"_FV2304" scheme (row 1, columns A:U).
$headers = @(
  "Segmentname_FV2210", "Segmentgruppe_FV2210", "Segment_FV2210", "Datenelement_FV2210", "Segment ID_FV2210",
  "Code_FV2210", "Qualifier_FV2210", "Beschreibung_FV2210", "Bedingungsausdruck_FV2210", "Bedingung_FV2210",
  "diff",
  "Segmentname_FV2304", "Segmentgruppe_FV2304", "Segment_FV2304", "Datenelement_FV2304", "Segment ID_FV2304",
  "Code_FV2304", "Qualifier_FV2304", "Beschreibung_FV2304", "Bedingungsausdruck_FV2304", "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Turn the used range into an Excel Table ("Table1") covering A1:U54,
# complete with header row + autofilter.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U54"), $null, 1)
$lo.Name = "Table1"

# Freeze the header row (split below row 1, keep it visible while scrolling).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
